$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.505.35"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "3.227.82"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("D5").Value = "'579.81"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").Value = "'182.22"
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.601"
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").Value = "3.225.74"
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("E10").Value = "  -2.98%  "
$ws.Range("E11").Value = "  -1.75%  "
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").Value = "3.787.66"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").Value = "'27.70"
$ws.Range("E15").Value = "  -3.33%  "
$ws.Range("D16").Value = "67.553.19"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").Value = "3.214.03"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("D20").Value = "'13.43"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("D21").Value = "'393.90"
$ws.Range("E21").Value = "  +3.07%  "
$ws.Range("E22").Value = "  -1.90%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'70.77"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").Value = "'0.513"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("E26").Value = "  -2.64%  "
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("E28").Value = "  -3.19%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").Value = "'1.96"
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("D31").Value = "'5.58"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").Value = "'22.62"
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("D33").Value = "'6.99"
$ws.Range("E33").Value = "  -4.13%  "
$ws.Range("E35").Value = "  -1.78%  "
$ws.Range("D36").Value = "'161.77"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("E37").Value = "  -5.31%  "
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").Value = "'26.30"
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("D40").Value = "'0.804"
$ws.Range("E40").Value = "  -3.72%  "
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("E42").Value = "  -4.57%  "
$ws.Range("E43").Value = "  -5.72%  "
$ws.Range("D44").Value = "'0.0681"
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("D46").Value = "2.607.71"
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("E47").Value = "  -3.47%  "
$ws.Range("D48").Value = "'334.56"
$ws.Range("E48").Value = "  -3.60%  "
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("E51").Value = "  -1.87%  "
